$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1 labels: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$oldNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$cols1 = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $ws.Range("$($cols1[$i])1").Value = "$($oldNames[$i])_FV2404"
}

$cols2 = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $ws.Range("$($cols2[$i])1").Value = "$($oldNames[$i])_FV2410"
}

# Freeze header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
